# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Update column G ("K") values for rows 2-38 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 0
    3  = 1
    4  = 1
    5  = 2
    6  = 2
    7  = 2
    8  = 1
    9  = 3
    10 = 1
    11 = 1
    12 = 0
    13 = 1
    14 = 0
    15 = 0
    16 = 2
    17 = 3
    18 = 0
    19 = 2
    20 = 4
    21 = 3
    22 = 1
    23 = 0
    24 = 1
    25 = 1
    26 = 0
    27 = 3
    28 = 0
    29 = 0
    30 = 1
    31 = 2
    32 = 0
    33 = 3
    34 = 1
    35 = 2
    36 = 1
    37 = 2
    38 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}

Write-Host "Updated column G (K) for rows 2-38"
